# Apply edits to the "Mentee" sheet as described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mentee")

# --- Update existing rows 32-34 ---
$ws.Range("A32").Value = "Trường hợp 1"
$ws.Range("C32").Value = "Kinh doanh quốc tế"

$ws.Range("A33").Value = "Trường hợp 2"
$ws.Range("B33").Value = "Nữ"
$ws.Range("C33").Value = "Tín dụng ngân hàng"

$ws.Range("A34").Value = "Trường hợp 3"
$ws.Range("C34").Value = "Web development"

# --- Add new rows 35-39 ---
$ws.Range("A35").Value = "Trường hợp 3"
$ws.Range("B35").Value = "Nam"
$ws.Range("C35").Value = "Web development"
$ws.Range("D35").Value = "Mentee"

$ws.Range("A36").Value = "Trường hợp 3"
$ws.Range("B36").Value = "Nam"
$ws.Range("C36").Value = "An toàn ứng dụng"
$ws.Range("D36").Value = "Mentee"

$ws.Range("A37").Value = "Trường hợp 4"
$ws.Range("B37").Value = "Nam"
$ws.Range("C37").Value = "Ngân hàng thương mại"
$ws.Range("D37").Value = "Mentee"

$ws.Range("A38").Value = "Trường hợp 5"
$ws.Range("B38").Value = "Nam"
$ws.Range("C38").Value = "Content creator"
$ws.Range("D38").Value = "Mentee"

$ws.Range("A39").Value = "Trường hợp 5"
$ws.Range("B39").Value = "Nữ"
$ws.Range("C39").Value = "Content creator"
$ws.Range("D39").Value = "Mentee"

# --- Window size change in workbook view (reflects the app window being
#     resized/maximized by the author before saving) ---
$excel.ActiveWindow.Width = 28800
$excel.ActiveWindow.Height = 12885
